# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
#  - inserts a new "Player Info" worksheet in front of the existing
#    "ODI Batting" / "ODI Bowling" sheets, with the player's basic info
#  - renames the MATCH_CARD_LINK column to MATCH_CODE on both the
#    "ODI Batting" and "ODI Bowling" sheets, and replaces the full
#    howstat.com scorecard URL with the bare match code it encoded.

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- new "Player Info" sheet, inserted before "ODI Batting" ---------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# leading "'" keeps the numeric-looking id stored as text, matching the
# scraped data's inline-string cells instead of turning it into a number
$playerInfo.Range("A2").Value = "'4494"
$playerInfo.Range("B2").Value = "Zafar Gohar"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# Worksheets.Add() shifts everybody's position, so re-resolve the other
# sheets fresh instead of reusing references captured before the insert
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# match the bold / bordered / centered header look already used by the
# other sheets, by copying the formatting from an existing header cell
$battingSheet.Range("A1").Copy() | Out-Null
$playerInfo.Range("A1:D1").PasteSpecial(-4122) | Out-Null

# --- ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ----------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").Value = "'3861"

# --- ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE ----------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").Value = "'3861"
